$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column L (12) - shifting
# intervention_type .. comment two columns to the right, so we can
# add "age_at_intervention" and "current_age" ahead of them.
$ws.Range("L1:M1").EntireColumn.Insert()

# New header values
$ws.Range("L1").Value = "age_at_intervention"
$ws.Range("M1").Value = "current_age"

# Match the final column widths from the template (engine quantizes
# ColumnWidth to a pixel grid, so feed the value that lands closest to
# the template's stored width).
$ws.Range("L1").EntireColumn.ColumnWidth = 16.92
$ws.Range("M1").EntireColumn.ColumnWidth = 14.92

# Reflect the active cell selection saved with the workbook
$ws.Range("A2").Select()
